$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, $Row, $Col, $Text)
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

Set-TextValue $ws 2 4 "275.36"
Set-TextValue $ws 2 5 "-1.41%"
Set-TextValue $ws 3 4 "27.32"
Set-TextValue $ws 3 5 "0.92%"
Set-TextValue $ws 4 4 "4.747"
Set-TextValue $ws 4 5 "-3.65%"
Set-TextValue $ws 5 4 "0.06316"
Set-TextValue $ws 5 5 "-0.80%"
Set-TextValue $ws 6 4 "6.927"
Set-TextValue $ws 6 5 "-0.23%"
Set-TextValue $ws 7 2 "GateToken"
Set-TextValue $ws 7 3 "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws 7 4 "3.302"
Set-TextValue $ws 7 5 "-1.73%"
Set-TextValue $ws 8 2 "FTXToken"
Set-TextValue $ws 8 3 "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws 8 4 "1.307"
Set-TextValue $ws 8 5 "38.28%"
Set-TextValue $ws 9 2 "MXToken"
Set-TextValue $ws 9 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws 9 4 "0.8782"
Set-TextValue $ws 9 5 "-0.83%"
Set-TextValue $ws 10 2 "WazirX"
Set-TextValue $ws 10 3 "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws 10 4 "0.1521"
Set-TextValue $ws 10 5 "3.55%"
Set-TextValue $ws 11 2 "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws 11 3 "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws 11 4 "0.04984"
Set-TextValue $ws 11 5 "-3.03%"
Set-TextValue $ws 12 2 "MandalaExchangeToken"
Set-TextValue $ws 12 3 "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws 12 4 "0.07475"
Set-TextValue $ws 12 5 "1.17%"
Set-TextValue $ws 13 2 "BitrueCoin"
Set-TextValue $ws 13 3 "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws 13 4 "0.02896"
Set-TextValue $ws 13 5 "-8.25%"
Set-TextValue $ws 14 2 "BitMartToken"
Set-TextValue $ws 14 3 "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws 14 4 "0.09055"
Set-TextValue $ws 14 5 "0.16%"
Set-TextValue $ws 15 2 "BitForexToken"
Set-TextValue $ws 15 3 "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws 15 4 "0.001568"
Set-TextValue $ws 15 5 "0.56%"
Set-TextValue $ws 16 2 "One"
Set-TextValue $ws 16 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws 16 4 "0.0006340"
Set-TextValue $ws 16 5 "1.10%"
Set-TextValue $ws 17 2 "TigerCash"
Set-TextValue $ws 17 3 "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws 17 4 "0.005936"
Set-TextValue $ws 17 5 "2.18%"
Set-TextValue $ws 18 2 "LEO"
Set-TextValue $ws 18 3 "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws 18 4 "3.452"
Set-TextValue $ws 18 5 "-1.06%"
Set-TextValue $ws 19 5 "-0.27%"
Set-TextValue $ws 20 5 "0.65%"
Set-TextValue $ws 21 4 "0.1326"
Set-TextValue $ws 21 5 "1.86%"
Set-TextValue $ws 22 4 "3.905"
Set-TextValue $ws 22 5 "0.52%"
Set-TextValue $ws 23 4 "0.04398"
Set-TextValue $ws 23 5 "1.13%"
Set-TextValue $ws 24 5 "-0.50%"
Set-TextValue $ws 25 4 "0.003840"
Set-TextValue $ws 25 5 "5.54%"
Set-TextValue $ws 26 4 "0.0001201"
Set-TextValue $ws 26 5 "0.15%"
Set-TextValue $ws 27 4 "0.0001937"
Set-TextValue $ws 27 5 "14.30%"
Set-TextValue $ws 40 4 "0.04112"
Set-TextValue $ws 40 5 "1.12%"
Set-TextValue $ws 41 4 "0.006873"
Set-TextValue $ws 41 5 "4.17%"
Set-TextValue $ws 42 4 "0.1174"
Set-TextValue $ws 42 5 "0.65%"
Set-TextValue $ws 43 4 "0.001932"
Set-TextValue $ws 43 5 "-17.75%"
Set-TextValue $ws 44 4 "0.01150"
Set-TextValue $ws 44 5 "-8.26%"
Set-TextValue $ws 45 4 "0.00005247"
Set-TextValue $ws 45 5 "-0.37%"
Set-TextValue $ws 46 2 "CoinbaseStockToken"
Set-TextValue $ws 46 3 "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue $ws 46 4 "0.02000"
Set-TextValue $ws 46 5 "-11.53%"
Set-TextValue $ws 47 2 "BOLO"
Set-TextValue $ws 47 3 "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue $ws 47 4 "1.490"
Set-TextValue $ws 47 5 "-37.23%"
